# Add columns I (I0) and J (IF) to the active worksheet, matching the
# style of the existing header row and filling in the data for rows 2-72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style of an existing header cell (e.g. H1) onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data values (rows 2-72) -------------------------------------------
$I = @(9,9,8,7,8,8,9,9,8,7,8,8,6,6,7,7,5,9,8,7,8,7,9,9,7,5,6,8,7,9,9,7,7,6,10,9,8,8,8,8,8,7,8,9,8,8,8,8,6,8,7,10,8,6,8,8,6,9,7,8,7,6,4,6,8,8,9,9,8,2,6)
$J = @(9,9,8,7,8,8,9,9,8,7,8,8,7,6,7,7,6,9,8,7,8,7,9,9,7,6,6,8,7,9,9,7,7,6,10,9,8,8,8,8,8,7,8,9,9,8,8,8,6,8,7,10,8,7,9,8,7,9,7,8,7,6,5,6,8,8,9,9,8,2,6)

for ($i = 0; $i -lt $I.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I[$i]
    $ws.Cells.Item($row, 10).Value = $J[$i]
}

# --- Update the sheet dimension to reflect the new used range ----------
$ws.Range("A1:J72").Select()
